$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 37 blank rows before row 5, pushing the existing rows 5-13 down
#    to rows 42-50 (matches the row relocation seen in the target workbook).
$ws.Rows("5:41").Insert()

# 2. The "6. Function ..." task entry (now at C50) is replaced with its
#    fuller description.
$ws.Range("C50").Value = '6. Функция която връща фраза, съдържаща думи'

# 3. Populate the freshly inserted row 5 with the new algorithm description;
#    reset its style to the workbook default first, then apply a red fill.
$ws.Range("C5").Style = "Normal"
$ws.Range("C5").Value = 'Информацията която се извлича за автор, факултетен номер и Университет се записва в структура (клас) със съответните полета и се замества в текста извлечен от първата страница на всеки документ с \n (символа за нов ред). След това в този текст се търсят ключови фрази "на тема", "тема на проекта", "тема на курсовия проект", "тема на курсовата работа", "тема на дипломната работа", "тема на преддипломния проект", "тема на преддипломната работа" и др И се извлича всичко след тях докато не се стигне до следните стоп фрази "\n\n" (!Независимо колко интервала има между тях - използвай невронната мрежа - PosTag), "по дисциплината", "разработен за" и др.'
$ws.Range("C5").Interior.Color = 255

# 4. The row insert pushed the two whole-column data validations down to
#    F1:F1042 / D1:D1042; restore them to their original extents.
$ws.Cells.Validation.Delete()

$fv = $ws.Range("F1:F1005").Validation
$fv.Add(3, 1, 3, '"Todo,In Progress,Testing,Done,"', 0)
$fv.IgnoreBlank = $false
$fv.ShowInput = $false
$fv.ShowError = $true

$dv = $ws.Range("D1:D1005").Validation
$dv.Add(3, 1, 3, '"Ванката,Ники"', 0)
$dv.IgnoreBlank = $false
$dv.ShowInput = $false
$dv.ShowError = $true

# 5. Restore the user selection to C11, matching the saved view state.
$ws.Range("C11").Select()
